$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the next day's gold-price row (row 14), mirroring the existing
# table rows: column A holds the new date, column B holds the (reused)
# price description text.

# Column A: "05-10-2025" looks like a date to Excel's input parser, so a
# plain Value assignment would get silently converted into a date serial
# number. Entering it as a formula that evaluates to a text string, then
# collapsing that formula down to its literal value via copy / paste-
# special (values only), keeps it as genuine text - matching how the
# existing date cells (A6:A13) are stored - without disturbing the
# cell's existing border formatting/style.
$ws.Range("A14").Formula = "=""05-10-2025"""
$ws.Range("A14").Copy()
$ws.Range("A14").PasteSpecial(-4163)

# Column B: plain descriptive text (not date-like), safe to assign
# directly; Excel will dedupe it against the identical existing shared
# string used by B13.
$ws.Range("B14").Value = "The price of gold in India today is ₹11,940 per gram for 24 karat gold, ₹10,945 per gram for 22 karat gold and ₹8,955 per gram for 18 karat gold (also called 999 gold)."
